$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing date style (A19 already carries the date number format)
# for the new date cells, instead of letting Excel mint a brand-new style.
$ws.Range("A19").Copy()
$ws.Range("A20:A21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 20: 12/10/2023, MPAL, TP, ..., x (col F)
$ws.Range("A20").Value = 45211
$ws.Range("B20").Value = "MPAL"
$ws.Range("C20").Value = "TP"
$ws.Range("F20").Value = "x"

# Row 21: 12/10/2023, MPAL, TP, x (col D)
$ws.Range("A21").Value = 45211
$ws.Range("B21").Value = "MPAL"
$ws.Range("C21").Value = "TP"
$ws.Range("D21").Value = "x"

# Update the view selection/top-left cell to match
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G21").Select()
